# Applies the "renewal of data w.r.t. current updates" edit:
#  - obj_output (sheet1) and rel_report__output (sheet3) each lose two
#    trailing rows (res_proc_costs / res_start_up_costs no longer needed,
#    ramp_costs takes the place previously held by res_proc_costs).
#  - "nonspin_units_starting_up" is renamed to "nonspin_units_started_up".
#  - Selections on all three sheets are refreshed.

$wb = $excel.ActiveWorkbook

$wsOutput = $wb.Worksheets.Item("obj_output")
$wsReport = $wb.Worksheets.Item("obj_report")
$wsRel    = $wb.Worksheets.Item("rel_report__output")

# --- obj_output (sheet1): drop the last two rows (res_proc_costs, res_start_up_costs) ---
$wsOutput.Rows.Item(33).Delete() | Out-Null
$wsOutput.Rows.Item(33).Delete() | Out-Null

# --- rel_report__output (sheet3): drop the last three rows (the two cost rows + blank row 35) ---
$wsRel.Rows.Item(35).Delete() | Out-Null
$wsRel.Rows.Item(34).Delete() | Out-Null
$wsRel.Rows.Item(33).Delete() | Out-Null

# --- rename remaining entries in place ---
$wsOutput.Range("B32").Value = "ramp_costs"
$wsRel.Range("C32").Value = "ramp_costs"

$wsOutput.Range("B19").Value = "nonspin_units_started_up"
$wsRel.Range("C19").Value = "nonspin_units_started_up"

# --- column width tweaks on rel_report__output (B/C got wider to fit the
#     longer strings now shown there after the row reshuffle) ---
$wsRel.Columns.Item(2).ColumnWidth = 33.833333333333336
$wsRel.Columns.Item(3).ColumnWidth = 27

# --- refresh selections ---
$wsOutput.Activate()
$wsOutput.Range("D11").Select() | Out-Null

$wsReport.Activate()
$wsReport.Range("A1:B2").Select() | Out-Null

$wsRel.Activate()
$wsRel.Range("F14").Select() | Out-Null
